# edit.ps1 — apply the "Removed vk link from front slide" commit.
#
# Changes applied (see diff):
#  1. Slide 1 ("front slide"): delete the "TextBox 7" shape that holds
#     the "http://vk.com/club33848893" link.
#  2. The three auto-updating "datetimeFigureOut" date placeholders
#     (slide master, slide layout 13, notes master) are re-stamped from
#     01.11.2012 to 06.11.2012 — this mirrors PowerPoint re-caching the
#     "Update automatically" date field text on save.

$p = $ppt.ActivePresentation

# --- 1. Remove the vk.com link textbox from the title/front slide -----
$frontSlide = $p.Slides.Item(1)
for ($i = $frontSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $frontSlide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 7") {
        $shp.Delete()
    }
}

# --- 2. Re-cache the "today" date field text on master/layout/notes ---
function Update-DatePlaceholder {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "01.11.2012") {
                $tr.Text = $newText
            }
        }
    }
}

# Slide master's own Date placeholder
Update-DatePlaceholder $p.SlideMaster.Shapes "06.11.2012"

# The "Занятие" custom layout (slideLayout13) Date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes "06.11.2012"
}

# Notes master Date placeholder
Update-DatePlaceholder $p.NotesMaster.Shapes "06.11.2012"
